$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update assessed value formula to use J2 directly (no longer divide by 100).
# Do this BEFORE re-typing J2 as text so M2 doesn't inherit a text style
# from its precedent.
$ws.Range("M2").Formula = "=L2*(J2)"

# Force text storage so the formatted strings aren't re-parsed as numbers
$ws.Range("G2:K2").NumberFormat = "@"

# Convert numeric value cells to formatted text strings
$ws.Range("G2").Value = "664,400.00"
$ws.Range("H2").Value = "93,300.00"
$ws.Range("I2").Value = "0.00"

# Simplify assessment rate display (drop decimal)
$ws.Range("J2").Value = "40%"

# Pad tax rate decimal places
$ws.Range("K2").Value = "3.25400"
